$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 415, pushing existing rows 415:488 down to 416:489.
$ws.Rows.Item(415).Insert()

# Populate the newly inserted row 415 with the new record.
$ws.Cells.Item(415, 1).Value = 3
$ws.Cells.Item(415, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(415, 3).Value = "Coquimbo"
$ws.Cells.Item(415, 4).Value = 44637
$ws.Cells.Item(415, 5).Value = 5
$ws.Cells.Item(415, 6).Value = 100112045
$ws.Cells.Item(415, 7).Value = "Zapallo"
$ws.Cells.Item(415, 8).Value = "Camote"
$ws.Cells.Item(415, 9).Value = "1a (cosecha)"
$ws.Cells.Item(415, 10).Value = 80
$ws.Cells.Item(415, 11).Value = 450
$ws.Cells.Item(415, 12).Value = 450
$ws.Cells.Item(415, 13).Value = 450
$ws.Cells.Item(415, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(415, 15).Value = "Provincia de Talca"
$ws.Cells.Item(415, 16).Value = 450
$ws.Cells.Item(415, 17).Value = 1
$ws.Cells.Item(415, 18).Value = "Hortaliza"
